# Applies the "Refactored the entire code" commit to SampleQuestions_Mistral.xlsx
# Changes:
#  - Header row (row 2): tweak wording of a few headers
#  - Keyword column (C): remove slashes (replace with spaces) in multi-word keywords,
#    drop the "/DNA" suffix from the nucleus keyword, and turn two keyword cells into
#    plain numbers instead of text
#  - Column F width narrowed
#  - Active cell/selection left on C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) ---
$ws.Range("B2").Value = "Expected Answer"
$ws.Range("D2").Value = "Similarity score %"
$ws.Range("F2").Value = "ActualAnswer"

# --- Keyword column (C) updates ---
$ws.Range("C11").Value = 3.14
$ws.Range("C13").Value = "JavaScript HTML CSS"
$ws.Range("C17").Value = "attention careful"
$ws.Range("C19").Value = "Rendering Navigation Interaction"
$ws.Range("C20").Value = "Cacti Cactuses"
$ws.Range("C21").Value = "nucleus"
$ws.Range("C22").Value = 299

# --- Column F width ---
$ws.Range("F:F").ColumnWidth = 14.2

# --- Selection ---
$ws.Range("C5").Select() | Out-Null
